$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.874.97"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "3.867.05"
$ws.Range("E3").Value = "  +3.14%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'602.93"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").Value = "'163.07"
$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("D7").Value = "3.863.18"
$ws.Range("E7").Value = "  +3.05%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").Value = "'6.31"
$ws.Range("E11").Value = "  -2.60%  "

$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").Value = "'36.86"
$ws.Range("E13").Value = "  -2.64%  "

$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").Value = "4.513.29"
$ws.Range("E15").Value = "  +3.12%  "

$ws.Range("D16").Value = "3.885.55"
$ws.Range("E16").Value = "  +3.61%  "

$ws.Range("D17").Value = "69.076.00"
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").Value = "'7.54"
$ws.Range("E18").Value = "  +2.66%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.113"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'11.42"
$ws.Range("E20").Value = "  +5.73%  "

$ws.Range("D21").Value = "'17.15"
$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("D22").Value = "'485.66"
$ws.Range("E22").Value = "  -1.25%  "

$ws.Range("D23").Value = "'0.721"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("E24").Value = "  +4.28%  "

$ws.Range("D25").Value = "'84.06"
$ws.Range("E25").Value = "  -0.80%  "

$ws.Range("E26").Value = "  -1.96%  "

$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").Value = "4.019.34"
$ws.Range("E32").Value = "  +3.19%  "

$ws.Range("E33").Value = "  -3.48%  "

$ws.Range("D34").Value = "'32.39"
$ws.Range("E34").Value = "  +2.69%  "

$ws.Range("D35").Value = "3.819.10"
$ws.Range("E35").Value = "  +3.64%  "

$ws.Range("E36").Value = "  -1.78%  "

$ws.Range("E37").Value = "  +1.77%  "

$ws.Range("E38").Value = "  +4.38%  "

$ws.Range("D39").Value = "'5.90"
$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").Value = "'443.15"
$ws.Range("E42").Value = "  +3.34%  "

$ws.Range("E43").Value = "  +1.33%  "

$ws.Range("D44").Value = "'48.54"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("E45").Value = "  -0.74%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "'8.40"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").Value = "'27.26"
$ws.Range("E48").Value = "  +15.55%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'143.18"
$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.833.79"
$ws.Range("E50").Value = "  +1.93%  "

$ws.Range("E51").Value = "  +1.22%  "
